$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Mon Sep 11 14:05:26 EDT 2023"
$ws.Range("B3").Value = "Mon Sep 11 14:05:40 EDT 2023"
$ws.Range("B4").Value = "Mon Sep 11 14:05:55 EDT 2023"
